$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header row + row labels to the new, more descriptive text.
# (Column B/C/D numeric data is untouched.)
$ws.Range("A1").Value = "Project(version)"
$ws.Range("B1").Value = "Total Statement Coverage"
$ws.Range("C1").Value = "Total Branch Coverage"
$ws.Range("D1").Value = "BMI"

$ws.Range("A2").Value  = "Apache commons Lang 3.0"
$ws.Range("A3").Value  = "Apache commons Lang 3.6"
$ws.Range("A4").Value  = "Apache commons Lang 3.7"
$ws.Range("A5").Value  = "Apache commons codec 1.11"
$ws.Range("A6").Value  = "Apache commons codec 1.12"
$ws.Range("A7").Value  = "Apache commons collections 4.0"
$ws.Range("A8").Value  = "Apache commons collections 4.4"
$ws.Range("A9").Value  = "Apache commons configuration 2.1"
$ws.Range("A10").Value = "Apache commons configuration 2.2"
$ws.Range("A11").Value = "Apache commons configuration 2.3"
$ws.Range("A12").Value = "Jfreechart 1.0.19"
$ws.Range("A13").Value = "jfreechart 1.5.0"

# Widen the columns to fit the new, longer labels.
# (values chosen so the engine's MDW-7 pixel-grid rounding of ColumnWidth
# lands as close as possible to the authored <col width="..."> targets of
# 35.6640625 / 31.5 / 31.33203125 / 11.83203125.)
$ws.Columns.Item(1).ColumnWidth = 35
$ws.Columns.Item(2).ColumnWidth = 30.857142857142858
$ws.Columns.Item(3).ColumnWidth = 30.571428571428573
$ws.Columns.Item(4).ColumnWidth = 11.142857142857142

# Move the active selection (matches the saved cursor position in the diff).
$ws.Range("C19").Select()
